$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new locatorType value for the "Hello, Sign in" wait step (E3, newly populated)
$ws.Range("E3").Value = "span"

# Change locatorType for the click step on row 4 from "div" to "span"
$ws.Range("E4").Value = "span"

# Set the new locatorType value for the "Continue" wait step (E5, newly populated)
$ws.Range("E5").Value = "span"

# Update the selected cell / scroll position to reflect the active work area
$ws.Activate()
$ws.Range("D7").Select()
